$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column F ("eta") between the existing eta_max (E) and
# pt_min (F, shifting to G) columns. Excel's column insert shifts all
# subsequent columns (old F:O -> new G:P) and adjusts formula refs,
# dimension, and the styled-column range automatically.
$ws.Columns("F").Insert()

# Header for the new column.
$ws.Range("F1").Value = "eta"

# Per-row eta values for the new column.
$etaValues = @(0.1, 0.3, 0.5, 0.71, 0.89, 1.11, 1.3, 1.49, 1.66, 1.88)
for ($i = 0; $i -lt $etaValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $etaValues[$i]
}

# The insert breaks the third shared-formula group (old M -> new N,
# "=K#/100"); restore it as a single range formula so Excel re-shares it
# the same way the original file had it shared across N2:N11.
$ws.Range("N2:N11").Formula = "=K2/100"

# Restore the active selection shown in the saved file.
[void]$ws.Range("I14").Select()
